# Auto-generated edit script: update crypto price/volume/hour data (GitHub Actions symbol-list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.63"
$ws.Range("E2").Value = "'-1.07%"
$ws.Range("G2").Value = "'7"
$ws.Range("E3").Value = "'-0.41%"
$ws.Range("G3").Value = "'7"
$ws.Range("D4").Value = "'5.040"
$ws.Range("E4").Value = "'-1.33%"
$ws.Range("G4").Value = "'7"
$ws.Range("D5").Value = "'0.07962"
$ws.Range("E5").Value = "'-1.63%"
$ws.Range("G5").Value = "'7"
$ws.Range("D6").Value = "'1.904"
$ws.Range("E6").Value = "'-2.73%"
$ws.Range("G6").Value = "'7"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'7.777"
$ws.Range("E7").Value = "'0.23%"
$ws.Range("G7").Value = "'7"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9211"
$ws.Range("E8").Value = "'-1.01%"
$ws.Range("G8").Value = "'7"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1336"
$ws.Range("E9").Value = "'-2.95%"
$ws.Range("G9").Value = "'7"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1903"
$ws.Range("E10").Value = "'-1.06%"
$ws.Range("G10").Value = "'7"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09129"
$ws.Range("E11").Value = "'-1.18%"
$ws.Range("G11").Value = "'7"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03442"
$ws.Range("E12").Value = "'0.37%"
$ws.Range("G12").Value = "'7"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09843"
$ws.Range("E13").Value = "'-0.02%"
$ws.Range("G13").Value = "'7"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001410"
$ws.Range("E14").Value = "'-2.06%"
$ws.Range("G14").Value = "'7"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006058"
$ws.Range("E15").Value = "'4.73%"
$ws.Range("G15").Value = "'7"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.731"
$ws.Range("E16").Value = "'3.02%"
$ws.Range("G16").Value = "'7"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.124"
$ws.Range("E17").Value = "'-1.48%"
$ws.Range("G17").Value = "'7"
$ws.Range("E18").Value = "'13.13%"
$ws.Range("G18").Value = "'7"
$ws.Range("E19").Value = "'0.06%"
$ws.Range("G19").Value = "'7"
$ws.Range("E20").Value = "'-2.33%"
$ws.Range("G20").Value = "'7"
$ws.Range("D21").Value = "'5.172"
$ws.Range("E21").Value = "'5.74%"
$ws.Range("G21").Value = "'7"
$ws.Range("D22").Value = "'0.2349"
$ws.Range("E22").Value = "'-6.09%"
$ws.Range("G22").Value = "'7"
$ws.Range("D23").Value = "'0.04418"
$ws.Range("E23").Value = "'-0.63%"
$ws.Range("G23").Value = "'7"
$ws.Range("D24").Value = "'0.001234"
$ws.Range("E24").Value = "'1.17%"
$ws.Range("G24").Value = "'7"
$ws.Range("D25").Value = "'0.004627"
$ws.Range("E25").Value = "'-4.33%"
$ws.Range("G25").Value = "'7"
$ws.Range("D26").Value = "'0.0001250"
$ws.Range("E26").Value = "'0.61%"
$ws.Range("G26").Value = "'7"
$ws.Range("D27").Value = "'0.0004440"
$ws.Range("E27").Value = "'-0.03%"
$ws.Range("G27").Value = "'7"
$ws.Range("G28").Value = "'7"
$ws.Range("G29").Value = "'7"
$ws.Range("G30").Value = "'7"
$ws.Range("G31").Value = "'7"
$ws.Range("G32").Value = "'7"
$ws.Range("G33").Value = "'7"
$ws.Range("G34").Value = "'7"
$ws.Range("G35").Value = "'7"
$ws.Range("G36").Value = "'7"
$ws.Range("G37").Value = "'7"
$ws.Range("G38").Value = "'7"
$ws.Range("D39").Value = "'0.01941"
$ws.Range("E39").Value = "'-4.02%"
$ws.Range("G39").Value = "'7"
$ws.Range("D40").Value = "'0.05379"
$ws.Range("E40").Value = "'8.97%"
$ws.Range("G40").Value = "'7"
$ws.Range("D41").Value = "'0.007589"
$ws.Range("E41").Value = "'-2.26%"
$ws.Range("G41").Value = "'7"
$ws.Range("D42").Value = "'0.01013"
$ws.Range("E42").Value = "'-0.88%"
$ws.Range("G42").Value = "'7"
$ws.Range("D43").Value = "'0.1354"
$ws.Range("E43").Value = "'-1.82%"
$ws.Range("G43").Value = "'7"
$ws.Range("D44").Value = "'0.002160"
$ws.Range("E44").Value = "'2.66%"
$ws.Range("G44").Value = "'7"
$ws.Range("D45").Value = "'0.01020"
$ws.Range("E45").Value = "'-11.78%"
$ws.Range("G45").Value = "'7"
$ws.Range("D46").Value = "'0.00006120"
$ws.Range("E46").Value = "'-5.02%"
$ws.Range("G46").Value = "'7"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.21%"
$ws.Range("G47").Value = "'7"
$ws.Range("G48").Value = "'7"
$ws.Range("D49").Value = "'0.001658"
$ws.Range("E49").Value = "'38.99%"
$ws.Range("G49").Value = "'7"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.21%"
$ws.Range("G50").Value = "'7"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.21%"
$ws.Range("G51").Value = "'7"
